# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted at N,
# pushing the old "Late" column (N) to O and the old "Outstanding" column
# (P) to Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N ("Late"), shifting
# N -> O and P -> Q.
$ws.Columns("N:N").Insert() | Out-Null

# Match the target column width for the newly inserted column N
# (raw OOXML width of 10, i.e. ColumnWidth = 10 - 5/6).
$ws.Columns("N:N").ColumnWidth = 9.166666666666666

# The selection after the edit is parked on R9.
$ws.Range("R9").Select() | Out-Null
